$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New interviewer "Niranjan" - several interview test rows recorded,
# only the first one (row 8) has the interview name filled in so far.
$ws.Range("A8").Value = "Niranjan"
$ws.Range("B8").Value = "1725184560-ITSE-Bad"

$ws.Range("A9").Value = "Niranjan"
$ws.Range("A10").Value = "Niranjan"
$ws.Range("A11").Value = "Niranjan"
$ws.Range("A12").Value = "Niranjan"
$ws.Range("A13").Value = "Niranjan"

# Reflect the last selected cell while entering this data
[void]$ws.Range("D8").Select()
